$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("B2").Value = 17.773194156146
$ws.Range("C2").Value = 7.791983378675781
$ws.Range("D2").Value = 8.239479812129757
$ws.Range("E2").Value = 12.36374021842432
$ws.Range("F2").Value = 33.83264864516157
$ws.Range("I2").Value = 24.73480789260679
$ws.Range("J2").Value = 9.581487980166377
$ws.Range("L2").Value = 11.48276104519299
$ws.Range("N2").Value = 18.03796723850547
$ws.Range("O2").Value = 25.9092512262303

# Row 3
$ws.Range("B3").Value = 17.34527673506544
$ws.Range("C3").Value = 7.480115152201104
$ws.Range("D3").Value = 8.239041773105118
$ws.Range("E3").Value = 12.39110689749546
$ws.Range("F3").Value = 33.88638093728091
$ws.Range("I3").Value = 24.83563599051509
$ws.Range("J3").Value = 9.601598045439774
$ws.Range("L3").Value = 11.46747575642974
$ws.Range("N3").Value = 18.08465470113758
$ws.Range("O3").Value = 25.97065680667476

# Row 4
$ws.Range("B4").Value = 17.07936528284242
$ws.Range("C4").Value = 7.280451675780372
$ws.Range("D4").Value = 8.239636584563442
$ws.Range("E4").Value = 12.40921478516403
$ws.Range("F4").Value = 33.92738678865055
$ws.Range("I4").Value = 24.90240761032847
$ws.Range("J4").Value = 9.614612442509172
$ws.Range("L4").Value = 11.45953316998233
$ws.Range("N4").Value = 18.11510521220662
$ws.Range("O4").Value = 26.01416566149579

# Row 5
$ws.Range("B5").Value = 16.97036124088553
$ws.Range("C5").Value = 7.197105834763712
$ws.Range("D5").Value = 8.240097022751005
$ws.Range("E5").Value = 12.4169224431989
$ws.Range("F5").Value = 33.94610871400051
$ws.Range("I5").Value = 24.93083900411933
$ws.Range("J5").Value = 9.620084030234183
$ws.Range("L5").Value = 11.45666184680697
$ws.Range("N5").Value = 18.12796355970592
$ws.Range("O5").Value = 26.03335261239151

# Row 6
$ws.Range("B6").Value = 16.95222699757571
$ws.Range("C6").Value = 7.183149023794884
$ws.Range("D6").Value = 8.240186671440535
$ws.Range("E6").Value = 12.41822215161759
$ws.Range("F6").Value = 33.94933886298966
$ws.Range("I6").Value = 24.93563374640724
$ws.Range("J6").Value = 9.621002750475995
$ws.Range("L6").Value = 11.45620720349895
$ws.Range("N6").Value = 18.13012585379857
$ws.Range("O6").Value = 26.03662647772468

# Row 7
$ws.Range("B7").Value = 17.07789761132716
$ws.Range("C7").Value = 7.279335561414481
$ws.Range("D7").Value = 8.239641910239058
$ws.Range("E7").Value = 12.40931740241285
$ws.Range("F7").Value = 33.92763113886452
$ws.Range("I7").Value = 24.9027861023848
$ws.Range("J7").Value = 9.614685552951471
$ws.Range("L7").Value = 11.45949296377202
$ws.Range("N7").Value = 18.11527680315373
$ws.Range("O7").Value = 26.01441852949608

# Row 8
$ws.Range("B8").Value = 17.62639587843995
$ws.Range("C8").Value = 7.686189513204181
$ws.Range("D8").Value = 8.239150021737963
$ws.Range("E8").Value = 12.37290576458991
$ws.Range("F8").Value = 33.84951033247263
$ws.Range("I8").Value = 24.76856333501077
$ws.Range("J8").Value = 9.588283826600232
$ws.Range("L8").Value = 11.47719285776912
$ws.Range("N8").Value = 18.05369527576094
$ws.Range("O8").Value = 25.9292171267847

# Row 9
$ws.Range("B9").Value = 18.67048084795372
$ws.Range("C9").Value = 8.416382288640627
$ws.Range("D9").Value = 8.245000668955406
$ws.Range("E9").Value = 12.31183269710247
$ws.Range("F9").Value = 33.76003277748693
$ws.Range("I9").Value = 24.54399768112118
$ws.Range("J9").Value = 9.541779362619828
$ws.Range("L9").Value = 11.52322913442204
$ws.Range("N9").Value = 17.94705256304993
$ws.Range("O9").Value = 25.80833122873813

# Row 10
$ws.Range("B10").Value = 19.41039411676416
$ws.Range("C10").Value = 8.908498719993192
$ws.Range("D10").Value = 8.253400555665902
$ws.Range("E10").Value = 12.27323016766589
$ws.Range("F10").Value = 33.73327545391878
$ws.Range("I10").Value = 24.40265162616972
$ws.Range("J10").Value = 9.510795320148732
$ws.Range("L10").Value = 11.56378698384841
$ws.Range("N10").Value = 17.8772563251918
$ws.Range("O10").Value = 25.74783805944563

# Row 11
$ws.Range("B11").Value = 19.73959102735308
$ws.Range("C11").Value = 9.122235878417753
$ws.Range("D11").Value = 8.258100060962267
$ws.Range("E11").Value = 12.25702360256396
$ws.Range("F11").Value = 33.72958175584942
$ws.Range("I11").Value = 24.34350166672237
$ws.Range("J11").Value = 9.497384728146637
$ws.Range("L11").Value = 11.58366025738705
$ws.Range("N11").Value = 17.84735075554639
$ws.Range("O11").Value = 25.72649582001936

# Row 12
$ws.Range("B12").Value = 19.86307287512799
$ws.Range("C12").Value = 9.201681942856572
$ws.Range("D12").Value = 8.260004786324712
$ws.Range("E12").Value = 12.25108080249981
$ws.Range("F12").Value = 33.72940198988066
$ws.Range("I12").Value = 24.32184518108778
$ws.Range("J12").Value = 9.492404404612541
$ws.Range("L12").Value = 11.59138670365392
$ws.Range("N12").Value = 17.83629082520811
$ws.Range("O12").Value = 25.71930390462433

# Row 13
$ws.Range("B13").Value = 19.8365329649705
$ws.Range("C13").Value = 9.184638648298328
$ws.Range("D13").Value = 8.259589024100281
$ws.Range("E13").Value = 12.25235205706769
$ws.Range("F13").Value = 33.7293865044983
$ws.Range("I13").Value = 24.32647624764379
$ws.Range("J13").Value = 9.493472655363002
$ws.Range("L13").Value = 11.58971380390681
$ws.Range("N13").Value = 17.83866101859551
$ws.Range("O13").Value = 25.72081320837611

# Row 14
$ws.Range("B14").Value = 19.74977418406484
$ws.Range("C14").Value = 9.128802028067714
$ws.Range("D14").Value = 8.25825426314667
$ws.Range("E14").Value = 12.25653079366875
$ws.Range("F14").Value = 33.72954254085744
$ws.Range("I14").Value = 24.34170508654647
$ws.Range("J14").Value = 9.49697303267555
$ws.Range("L14").Value = 11.58429191345809
$ws.Range("N14").Value = 17.8464355479618
$ws.Range("O14").Value = 25.7258862925515

# Row 15
$ws.Range("B15").Value = 19.69647528380903
$ws.Range("C15").Value = 9.09440529588402
$ws.Range("D15").Value = 8.257452942421773
$ws.Range("E15").Value = 12.25911567836867
$ws.Range("F15").Value = 33.72979684130995
$ws.Range("I15").Value = 24.35112991839642
$ws.Range("J15").Value = 9.499129862164841
$ws.Range("L15").Value = 11.58099689224013
$ws.Range("N15").Value = 17.85123211891652
$ws.Range("O15").Value = 25.72910964396991

# Row 16
$ws.Range("B16").Value = 19.38872153362166
$ws.Range("C16").Value = 8.894323550732866
$ws.Range("D16").Value = 8.253111006966225
$ws.Range("E16").Value = 12.27431651343315
$ws.Range("F16").Value = 33.73368745186028
$ws.Range("I16").Value = 24.40662098891095
$ws.Range("J16").Value = 9.511685466576433
$ws.Range("L16").Value = 11.56251653840216
$ws.Range("N16").Value = 17.8792477953782
$ws.Range("O16").Value = 25.74935726972843

# Row 17
$ws.Range("B17").Value = 19.19794640109044
$ws.Range("C17").Value = 8.768958874030771
$ws.Range("D17").Value = 8.250671515118052
$ws.Range("E17").Value = 12.28398820746188
$ws.Range("F17").Value = 33.73824582807321
$ws.Range("I17").Value = 24.44198310502339
$ws.Range("J17").Value = 9.519562870518348
$ws.Range("L17").Value = 11.55154133511651
$ws.Range("N17").Value = 17.8969065997555
$ws.Range("O17").Value = 25.7633618343367

# Row 18
$ws.Range("B18").Value = 19.08752796346511
$ws.Range("C18").Value = 8.695900743450229
$ws.Range("D18").Value = 8.249351133133457
$ws.Range("E18").Value = 12.28967857230645
$ws.Range("F18").Value = 33.74166577204157
$ws.Range("I18").Value = 24.46280699988176
$ws.Range("J18").Value = 9.524158169785652
$ws.Range("L18").Value = 11.54536289552824
$ws.Range("N18").Value = 17.90723718373091
$ws.Range("O18").Value = 25.77199815969126

# Row 19
$ws.Range("B19").Value = 19.05002738069454
$ws.Range("C19").Value = 8.67100222197752
$ws.Range("D19").Value = 8.24891832180796
$ws.Range("E19").Value = 12.2916271365527
$ws.Range("F19").Value = 33.7429607716459
$ws.Range("I19").Value = 24.46994077586876
$ws.Range("J19").Value = 9.525725136788378
$ws.Range("L19").Value = 11.54329414510557
$ws.Range("N19").Value = 17.9107647967697
$ws.Range("O19").Value = 25.7750220448337

# Row 20
$ws.Range("B20").Value = 19.21832697375975
$ws.Range("C20").Value = 8.782402924843295
$ws.Range("D20").Value = 8.250922647344547
$ws.Range("E20").Value = 12.28294545028604
$ws.Range("F20").Value = 33.73767798112758
$ws.Range("I20").Value = 24.43816858554848
$ws.Range("J20").Value = 9.518717642679539
$ws.Range("L20").Value = 11.55269580128215
$ws.Range("N20").Value = 17.89500881673666
$ws.Range("O20").Value = 25.76181085122946

# Row 21
$ws.Range("B21").Value = 19.77529017692303
$ws.Range("C21").Value = 9.145243329244634
$ws.Range("D21").Value = 8.258642928642596
$ws.Range("E21").Value = 12.25529812821105
$ws.Range("F21").Value = 33.72946363287274
$ws.Range("I21").Value = 24.33721184592989
$ws.Range("J21").Value = 9.495942231207913
$ws.Range("L21").Value = 11.58587903375433
$ws.Range("N21").Value = 17.84414480260944
$ws.Range("O21").Value = 25.72437204015517

# Row 22
$ws.Range("B22").Value = 20.13238550825274
$ws.Range("C22").Value = 9.373674055645516
$ws.Range("D22").Value = 8.264417313428361
$ws.Range("E22").Value = 12.23836119119736
$ws.Range("F22").Value = 33.7311994113404
$ws.Range("I22").Value = 24.27555847550035
$ws.Range("J22").Value = 9.481628077973125
$ws.Range("L22").Value = 11.6087351254181
$ws.Range("N22").Value = 17.81244457652703
$ws.Range("O22").Value = 25.7050913435646

# Row 23
$ws.Range("B23").Value = 19.94246501385165
$ws.Range("C23").Value = 9.252563098081163
$ws.Range("D23").Value = 8.261269141227032
$ws.Range("E23").Value = 12.24729729605762
$ws.Range("F23").Value = 33.72962323289714
$ws.Range("I23").Value = 24.30806742801804
$ws.Range("J23").Value = 9.489215708822121
$ws.Range("L23").Value = 11.59643074289298
$ws.Range("N23").Value = 17.82922267225659
$ws.Range("O23").Value = 25.7149066502062

# Row 24
$ws.Range("B24").Value = 19.20911521068522
$ws.Range("C24").Value = 8.776327932704849
$ws.Range("D24").Value = 8.250808854610938
$ws.Range("E24").Value = 12.28341647590803
$ws.Range("F24").Value = 33.73793221502814
$ws.Range("I24").Value = 24.43989159194389
$ws.Range("J24").Value = 9.519099563090203
$ws.Range("L24").Value = 11.55217345770135
$ws.Range("N24").Value = 17.89586624897873
$ws.Range("O24").Value = 25.76251022883805

# Row 25
$ws.Range("B25").Value = 18.39225837097892
$ws.Range("C25").Value = 8.22644019292348
$ws.Range("D25").Value = 8.242693582126648
$ws.Range("E25").Value = 12.32725174165337
$ws.Range("F25").Value = 33.77740024798908
$ws.Range("I25").Value = 24.60060289914231
$ws.Range("J25").Value = 9.55379900650971
$ws.Range("L25").Value = 11.50957929675196
$ws.Range("N25").Value = 17.9743963012453
$ws.Range("O25").Value = 25.83607082150789
